$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 130
$ws.Range('B130').Value = 7483081
$ws.Range('F130').Value = 'Deportivo Cuenca'
$ws.Range('G130').Value = 'El Nacional'
$ws.Range('H130').Value = 1
$ws.Range('J130').Value = 'H'
$ws.Range('K130').Value = 2.75
$ws.Range('M130').Value = 2.55
$ws.Range('N130').Value = 3
$ws.Range('O130').Value = 3.3
$ws.Range('P130').Value = 2.3
$ws.Range('Q130').Value = 0.25
$ws.Range('R130').Value = 1.825
$ws.Range('S130').Value = 1.975
$ws.Range('T130').Value = 2.75
$ws.Range('U130').Value = 2
$ws.Range('V130').Value = 1.8
$ws.Range('W130').Value = 2
$ws.Range('X130').Value = -1
$ws.Range('Z130').Value = 0.825
$ws.Range('AA130').Value = -1
$ws.Range('AC130').Value = 0.8
# Row 132
$ws.Range('B132').Value = 7483247
$ws.Range('F132').Value = 'Mushuc Runa'
$ws.Range('G132').Value = 'Universidad Catolica del Ecuador'
$ws.Range('H132').Value = 0
$ws.Range('I132').Value = 2
$ws.Range('J132').Value = 'A'
$ws.Range('K132').Value = 3.25
$ws.Range('L132').Value = 3.2
$ws.Range('M132').Value = 2.25
$ws.Range('N132').Value = 3.5
$ws.Range('O132').Value = 3.25
$ws.Range('P132').Value = 2.1
$ws.Range('Q132').Value = 0.5
$ws.Range('R132').Value = 1.775
$ws.Range('S132').Value = 2.025
$ws.Range('T132').Value = 2.5
$ws.Range('U132').Value = 1.9
$ws.Range('V132').Value = 1.9
$ws.Range('W132').Value = -1
$ws.Range('Y132').Value = 1.1
$ws.Range('Z132').Value = -1
$ws.Range('AA132').Value = 1.025
$ws.Range('AC132').Value = 0.8999999999999999
# Row 133
$ws.Range('B133').Value = 7483281
$ws.Range('F133').Value = 'SD Aucas'
$ws.Range('G133').Value = 'Delfin SC'
$ws.Range('I133').Value = 0
$ws.Range('J133').Value = 'D'
$ws.Range('K133').Value = 1.909
$ws.Range('L133').Value = 3.25
$ws.Range('M133').Value = 4.2
$ws.Range('N133').Value = 1.909
$ws.Range('O133').Value = 3.5
$ws.Range('P133').Value = 4
$ws.Range('Q133').Value = -0.5
$ws.Range('R133').Value = 1.9
$ws.Range('S133').Value = 1.9
$ws.Range('U133').Value = 1.8
$ws.Range('V133').Value = 2
$ws.Range('X133').Value = 2.5
$ws.Range('Y133').Value = -1
$ws.Range('AA133').Value = 0.8999999999999999
$ws.Range('AC133').Value = 1
# Row 135
$ws.Range('B135').Value = 7482867
$ws.Range('F135').Value = 'Cumbaya FC'
$ws.Range('G135').Value = 'LDU Quito'
$ws.Range('I135').Value = 2
$ws.Range('J135').Value = 'A'
$ws.Range('K135').Value = 5.25
$ws.Range('L135').Value = 3.75
$ws.Range('M135').Value = 1.65
$ws.Range('N135').Value = 9
$ws.Range('O135').Value = 4.5
$ws.Range('P135').Value = 1.363
$ws.Range('Q135').Value = 1.25
$ws.Range('R135').Value = 1.975
$ws.Range('S135').Value = 1.825
$ws.Range('T135').Value = 2.5
$ws.Range('U135').Value = 1.825
$ws.Range('V135').Value = 1.975
$ws.Range('X135').Value = -1
$ws.Range('Y135').Value = 0.363
$ws.Range('Z135').Value = 0.4875
$ws.Range('AA135').Value = -0.5
$ws.Range('AB135').Value = 0.825
$ws.Range('AC135').Value = -1
# Row 137
$ws.Range('B137').Value = 7483306
$ws.Range('F137').Value = 'Tecnico Universitario'
$ws.Range('G137').Value = 'Club Atletico Libertad'
$ws.Range('I137').Value = 1
$ws.Range('J137').Value = 'D'
$ws.Range('K137').Value = 1.5
$ws.Range('L137').Value = 4.333
$ws.Range('M137').Value = 5.75
$ws.Range('N137').Value = 1.533
$ws.Range('O137').Value = 4.2
$ws.Range('P137').Value = 5.5
$ws.Range('Q137').Value = -1
$ws.Range('R137').Value = 1.925
$ws.Range('S137').Value = 1.875
$ws.Range('T137').Value = 2.25
$ws.Range('U137').Value = 1.8
$ws.Range('V137').Value = 2
$ws.Range('X137').Value = 3.2
$ws.Range('Y137').Value = -1
$ws.Range('Z137').Value = -1
$ws.Range('AA137').Value = 0.875
$ws.Range('AB137').Value = -0.5
$ws.Range('AC137').Value = 0.5
# Row 142
$ws.Range('B142').Value = 7528848
$ws.Range('F142').Value = 'Emelec'
$ws.Range('G142').Value = 'Deportivo Cuenca'
$ws.Range('H142').Value = 2
$ws.Range('I142').Value = 1
$ws.Range('J142').Value = 'H'
$ws.Range('K142').Value = 1.75
$ws.Range('L142').Value = 3.5
$ws.Range('M142').Value = 4.2
$ws.Range('N142').Value = 2.4
$ws.Range('O142').Value = 3.1
$ws.Range('P142').Value = 2.75
$ws.Range('Q142').Value = -0.25
$ws.Range('R142').Value = 2.05
$ws.Range('S142').Value = 1.75
$ws.Range('U142').Value = 1.8
$ws.Range('V142').Value = 2
$ws.Range('W142').Value = 1.4
$ws.Range('Y142').Value = -1
$ws.Range('Z142').Value = 1.05
$ws.Range('AA142').Value = -1
$ws.Range('AB142').Value = 0.8
# Row 143
$ws.Range('B143').Value = 7528852
$ws.Range('F143').Value = 'Delfin SC'
$ws.Range('G143').Value = 'Tecnico Universitario'
$ws.Range('H143').Value = 2
$ws.Range('I143').Value = 2
$ws.Range('J143').Value = 'D'
$ws.Range('K143').Value = 2.1
$ws.Range('L143').Value = 3.4
$ws.Range('M143').Value = 3.1
$ws.Range('N143').Value = 2.1
$ws.Range('O143').Value = 3.4
$ws.Range('P143').Value = 3.1
$ws.Range('Q143').Value = -0.25
$ws.Range('T143').Value = 2.25
$ws.Range('U143').Value = 1.9
$ws.Range('V143').Value = 1.9
$ws.Range('X143').Value = 2.4
$ws.Range('Y143').Value = -1
$ws.Range('Z143').Value = -0.5
$ws.Range('AA143').Value = 0.5
$ws.Range('AB143').Value = 0.8999999999999999
$ws.Range('AC143').Value = -1
# Row 144
$ws.Range('B144').Value = 7528857
$ws.Range('F144').Value = 'Universidad Catolica del Ecuador'
$ws.Range('G144').Value = 'Barcelona Guayaquil'
$ws.Range('H144').Value = 0
$ws.Range('I144').Value = 1
$ws.Range('J144').Value = 'A'
$ws.Range('K144').Value = 1.533
$ws.Range('L144').Value = 4
$ws.Range('M144').Value = 5.5
$ws.Range('N144').Value = 1.5
$ws.Range('O144').Value = 4.333
$ws.Range('P144').Value = 5.25
$ws.Range('Q144').Value = -1
$ws.Range('T144').Value = 3
$ws.Range('U144').Value = 1.975
$ws.Range('V144').Value = 1.825
$ws.Range('X144').Value = -1
$ws.Range('Y144').Value = 4.25
$ws.Range('Z144').Value = -1
$ws.Range('AA144').Value = 1
$ws.Range('AB144').Value = -1
$ws.Range('AC144').Value = 0.825
# Row 145
$ws.Range('B145').Value = 7528858
$ws.Range('F145').Value = 'Orense'
$ws.Range('G145').Value = 'SD Aucas'
$ws.Range('H145').Value = 1
$ws.Range('I145').Value = 2
$ws.Range('J145').Value = 'A'
$ws.Range('K145').Value = 2.2
$ws.Range('L145').Value = 3.2
$ws.Range('M145').Value = 3.2
$ws.Range('N145').Value = 1.95
$ws.Range('O145').Value = 3.2
$ws.Range('P145').Value = 3.8
$ws.Range('Q145').Value = -0.5
$ws.Range('R145').Value = 1.95
$ws.Range('S145').Value = 1.85
$ws.Range('U145').Value = 1.85
$ws.Range('V145').Value = 1.95
$ws.Range('W145').Value = -1
$ws.Range('Y145').Value = 2.8
$ws.Range('Z145').Value = -1
$ws.Range('AA145').Value = 0.8500000000000001
$ws.Range('AB145').Value = 0.8500000000000001
# Row 200
$ws.Range('U200').Value = 2
$ws.Range('V200').Value = 1.8
# Row 202
$ws.Range('N202').Value = 1.571
$ws.Range('P202').Value = 5
$ws.Range('Q202').Value = -0.75
$ws.Range('R202').Value = 1.775
$ws.Range('S202').Value = 2.025
# Row 204
$ws.Range('R204').Value = 1.85
$ws.Range('S204').Value = 1.95
# Row 206
$ws.Range('R206').Value = 1.8
$ws.Range('S206').Value = 2
